$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.121.72"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.28"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.89"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4674"
$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2827"
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06442"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.72"
$ws.Range("E10").Value = "  -4.31%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07667"
$ws.Range("E11").Value = "  -4.35%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.32"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.29"
$ws.Range("E13").Value = "  -3.84%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.031"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6758"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "266.52"
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.103.17"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.27"
$ws.Range("E18").Value = "  -5.23%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007494"
$ws.Range("E20").Value = "  -1.66%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.093.89"
$ws.Range("E21").Value = "  -0.77%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.128"
$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.063"
$ws.Range("E24").Value = "  -2.49%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.237"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.19"
$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.52"
$ws.Range("E27").Value = "  -2.07%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.871"
$ws.Range("E28").Value = "  -4.01%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("E29").Value = "  +0.44%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09771"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.446"
$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.175"
$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.952"
$ws.Range("E33").Value = "  -2.94%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04625"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.104"
$ws.Range("E35").Value = "  -3.14%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6794"
$ws.Range("E36").Value = "  -3.02%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01809"
$ws.Range("E38").Value = "  -3.22%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.706"
$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.270"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.04"
$ws.Range("E41").Value = "  -2.70%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8261"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.90"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.867"
$ws.Range("E45").Value = "  -4.97%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4021"
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.170"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "919.88"
$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.862"
$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.99"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05538"
$ws.Range("E51").Value = "  -2.78%  "
